$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data to append after the existing last row (row 22 -> new row 23)
$newRow = 23

$ws.Cells.Item($newRow, 1).Value = "GFG"
$ws.Cells.Item($newRow, 2).Value = "Segregate even and odd nodes in a Link List "

# Match the formatting/style used by the rest of the table rows (left aligned, wrap text)
# which is the same style already applied to both columns on row 16 (style index 6 in the XML).
$srcStyleRange = $ws.Range("A16:B16")
$dstStyleRange = $ws.Range("A$newRow`:B$newRow")
$srcStyleRange.Copy() | Out-Null
$dstStyleRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial formats only, but make sure values remain correct)
$ws.Cells.Item($newRow, 1).Value = "GFG"
$ws.Cells.Item($newRow, 2).Value = "Segregate even and odd nodes in a Link List "

# Update the visible selection to match the recorded state after the edit
$ws.Range("B19").Select() | Out-Null
